# Refresh scraped "want to go" counters (column F) across sheets,
# matching the upstream data re-generation captured in the diff.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 917
$wsExhibit.Range("F5").Value = 954
$wsExhibit.Range("F6").Value = 1708
$wsExhibit.Range("F7").Value = 372
$wsExhibit.Range("F16").Value = 128
$wsExhibit.Range("F18").Value = 23
$wsExhibit.Range("F29").Value = 96
$wsExhibit.Range("F31").Value = 245

# 演出 (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F11").Value = 111

# 全部类型 (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 917
$wsAll.Range("F6").Value = 954
$wsAll.Range("F7").Value = 1708
$wsAll.Range("F8").Value = 372
$wsAll.Range("F18").Value = 128
$wsAll.Range("F21").Value = 23
$wsAll.Range("F39").Value = 96
$wsAll.Range("F41").Value = 245
$wsAll.Range("F43").Value = 111
$wsAll.Range("F44").Value = 111
